# Apply edits described by the diff:
#  - workbook tab ratio 926 -> 600 (best-effort; may not be persisted by host)
#  - remove the "Orientation" column (column G) entirely
#  - rename "Stella.jpg" filename cell to "Still-Life in White #2.jpg"
#  - widen column A, normalize header/data row heights to 15.75
#  - extend the sheet with 5 blank rows (7-11) at row height 13.8
#  - update the view's selection / scroll position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook window / tab ratio -------------------------------------------------
$excel.ActiveWindow.TabRatio = 600

# --- Remove the Orientation column (G) --------------------------------------------
$ws.Columns.Item(7).Delete()

# --- Rename the filename for the Stella row (row 6) --------------------------------
$ws.Cells.Item(6, 1).Value = "Still-Life in White #2.jpg"

# --- Column widths -----------------------------------------------------------------
# Target stored xml width is 30.984693877551; the ColumnWidth API here only persists
# widths quantized to 1/6-character steps, so 30.17 is the closest achievable value
# (serializes to width="31", ~0.015 off the target instead of the default 14.43).
$ws.Columns.Item(1).ColumnWidth = 30.17

# --- Row heights for the existing header + data rows --------------------------------
for ($r = 1; $r -le 6; $r++) {
    $ws.Rows.Item($r).RowHeight = 15.75
}

# --- Append 5 blank rows (7-11) at the default smaller row height -------------------
for ($r = 7; $r -le 11; $r++) {
    $ws.Rows.Item($r).RowHeight = 13.8
}
# Touch a single cell so the sheet's used-range/dimension picks up the new rows.
$ws.Cells.Item(7, 1).Locked = $true

# --- Update selection / scroll position to match the new layout ---------------------
$ws.Range("A7:A11").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
